# Regenerate handback report: refresh the "076ce8df-ba51-40e2-b7bc-b5fcb9750b91.md"
# row's handoff/handback timestamps across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G, row 2) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-11-09 18:55:10"

# --- zh-cn sheet: Correspond Handoff / Handback DateTime (columns H & K, row 2) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-11-09 18:54:57"
$wsZhCn.Range("K2").Value = "2016-11-09 18:55:49"

# --- de-de sheet: Correspond Handoff / Handback DateTime (columns H & K, row 2) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-11-09 18:55:10"
$wsDeDe.Range("K2").Value = "2016-11-09 18:56:07"
